$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2, shifting existing rows (2:168) down to (3:169).
$ws.Rows("2:2").Insert()

# The row that used to be row 2 is now row 3. Duplicate that row's content (and
# formatting/hyperlink) up into the freshly-inserted row 2, then update the
# date in column A to the new latest date. This mirrors the source data feed,
# which always prepends "yesterday's newest row" content and bumps the date.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial()
$ws.Range("A2").Value = "26-11-2025"
